$d = $word.ActiveDocument

# Locate the (unique, case-sensitive) "Present" in the "12th: Everyone Present"
# attendance line and fix its capitalisation to "present".
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Present", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng = $find.Parent

    # Replace just the leading "P" with a lower-case "p". Toggling a
    # character-level formatting property on the single-letter range and
    # then reverting it forces Word to split the run at that character
    # boundary (matching how a manual retype of a single letter behaves),
    # producing the same "...Everyone " / "p" / "resent" run structure
    # seen when this fix was made by hand.
    $pRange = $d.Range($rng.Start, $rng.Start + 1)
    $pRange.Font.Bold = 1
    $pRange.Text = "p"

    $pRange2 = $d.Range($rng.Start, $rng.Start + 1)
    $pRange2.Font.Bold = 0
}
